# Re-applies the refreshed word/image/category cue list onto sheet1.
# Only the rows whose word, image path or category actually changed
# are touched below; everything else (headers, "none" filler rows,
# formatting) is left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 3: 'lügen'/'dog/dog010.jpg'/'dog' -> 'schaden'/'house/house009.jpg'/'house'
$ws.Range("A3").Value = "schaden"
$ws.Range("B3").Value = "house/house009.jpg"
$ws.Range("C3").Value = "house"

# row 4: 'wüten'/'house/house011.jpg'/'house' -> 'brauchen'/'dog/dog012.jpg'/'dog'
$ws.Range("A4").Value = "brauchen"
$ws.Range("B4").Value = "dog/dog012.jpg"
$ws.Range("C4").Value = "dog"

# row 6: 'gleichen'/'dog/dog006.jpg'/'dog' -> 'wohnen'/'dog/dog009.jpg'/'dog'
$ws.Range("A6").Value = "wohnen"
$ws.Range("B6").Value = "dog/dog009.jpg"

# row 7: 'ändern'/'dog/dog013.jpg'/'dog' -> 'packen'/'house/house026.jpg'/'house'
$ws.Range("A7").Value = "packen"
$ws.Range("B7").Value = "house/house026.jpg"
$ws.Range("C7").Value = "house"

# row 9: 'posten'/'dog/dog007.jpg'/'dog' -> 'albern'/'dog/dog011.jpg'/'dog'
$ws.Range("A9").Value = "albern"
$ws.Range("B9").Value = "dog/dog011.jpg"

# row 10: 'schweben'/'dog/dog024.jpg'/'dog' -> 'kleben'/'house/house021.jpg'/'house'
$ws.Range("A10").Value = "kleben"
$ws.Range("B10").Value = "house/house021.jpg"
$ws.Range("C10").Value = "house"

# row 12: 'wohnen'/'dog/dog001.jpg'/'dog' -> 'ändern'/'dog/dog021.jpg'/'dog'
$ws.Range("A12").Value = "ändern"
$ws.Range("B12").Value = "dog/dog021.jpg"

# row 13: 'beißen'/'house/house009.jpg'/'house' -> 'steuern'/'dog/dog020.jpg'/'dog'
$ws.Range("A13").Value = "steuern"
$ws.Range("B13").Value = "dog/dog020.jpg"
$ws.Range("C13").Value = "dog"

# row 15: 'sammeln'/'house/house023.jpg'/'house' -> 'parken'/'house/house003.jpg'/'house'
$ws.Range("A15").Value = "parken"
$ws.Range("B15").Value = "house/house003.jpg"

# row 16: 'schaden'/'house/house010.jpg'/'house' -> 'stürmen'/'house/house005.jpg'/'house'
$ws.Range("A16").Value = "stürmen"
$ws.Range("B16").Value = "house/house005.jpg"

# row 18: 'schreiben'/'house/house026.jpg'/'house' -> 'decken'/'dog/dog019.jpg'/'dog'
$ws.Range("A18").Value = "decken"
$ws.Range("B18").Value = "dog/dog019.jpg"
$ws.Range("C18").Value = "dog"

# row 19: 'fordern'/'house/house022.jpg'/'house' -> 'süßen'/'dog/dog008.jpg'/'dog'
$ws.Range("A19").Value = "süßen"
$ws.Range("B19").Value = "dog/dog008.jpg"
$ws.Range("C19").Value = "dog"

# row 21: 'süßen'/'dog/dog027.jpg'/'dog' -> 'landen'/'dog/dog022.jpg'/'dog'
$ws.Range("A21").Value = "landen"
$ws.Range("B21").Value = "dog/dog022.jpg"

# row 22: 'ruhen'/'dog/dog031.jpg'/'dog' -> 'triefen'/'dog/dog024.jpg'/'dog'
$ws.Range("A22").Value = "triefen"
$ws.Range("B22").Value = "dog/dog024.jpg"

# row 24: 'münzen'/'house/house029.jpg'/'house' -> 'rasen'/'dog/dog030.jpg'/'dog'
$ws.Range("A24").Value = "rasen"
$ws.Range("B24").Value = "dog/dog030.jpg"
$ws.Range("C24").Value = "dog"

# row 25: 'schützen'/'house/house019.jpg'/'house' -> 'fügen'/'dog/dog027.jpg'/'dog'
$ws.Range("A25").Value = "fügen"
$ws.Range("B25").Value = "dog/dog027.jpg"
$ws.Range("C25").Value = "dog"

# row 27: 'kleben'/'house/house018.jpg'/'house' -> 'bilden'/'house/house018.jpg'/'house'
$ws.Range("A27").Value = "bilden"

# row 28: 'nerven'/'dog/dog025.jpg'/'dog' -> 'werden'/'dog/dog010.jpg'/'dog'
$ws.Range("A28").Value = "werden"
$ws.Range("B28").Value = "dog/dog010.jpg"

# row 30: 'albern'/'dog/dog005.jpg'/'dog' -> 'klingen'/'house/house012.jpg'/'house'
$ws.Range("A30").Value = "klingen"
$ws.Range("B30").Value = "house/house012.jpg"
$ws.Range("C30").Value = "house"

# row 31: 'bitten'/'house/house003.jpg'/'house' -> 'schreiben'/'house/house001.jpg'/'house'
$ws.Range("A31").Value = "schreiben"
$ws.Range("B31").Value = "house/house001.jpg"

# row 33: 'stillen'/'dog/dog004.jpg'/'dog' -> 'starren'/'house/house004.jpg'/'house'
$ws.Range("A33").Value = "starren"
$ws.Range("B33").Value = "house/house004.jpg"
$ws.Range("C33").Value = "house"

# row 34: 'deuten'/'dog/dog021.jpg'/'dog' -> 'lassen'/'house/house019.jpg'/'house'
$ws.Range("A34").Value = "lassen"
$ws.Range("B34").Value = "house/house019.jpg"
$ws.Range("C34").Value = "house"

# row 36: 'biegen'/'dog/dog022.jpg'/'dog' -> 'fordern'/'house/house023.jpg'/'house'
$ws.Range("A36").Value = "fordern"
$ws.Range("B36").Value = "house/house023.jpg"
$ws.Range("C36").Value = "house"

# row 37: 'tauschen'/'dog/dog011.jpg'/'dog' -> 'stillen'/'dog/dog023.jpg'/'dog'
$ws.Range("A37").Value = "stillen"
$ws.Range("B37").Value = "dog/dog023.jpg"

# row 39: 'triefen'/'dog/dog002.jpg'/'dog' -> 'herrschen'/'house/house016.jpg'/'house'
$ws.Range("A39").Value = "herrschen"
$ws.Range("B39").Value = "house/house016.jpg"
$ws.Range("C39").Value = "house"

# row 40: 'landen'/'dog/dog029.jpg'/'dog' -> 'sammeln'/'house/house017.jpg'/'house'
$ws.Range("A40").Value = "sammeln"
$ws.Range("B40").Value = "house/house017.jpg"
$ws.Range("C40").Value = "house"

# row 42: 'lassen'/'house/house014.jpg'/'house' -> 'schultern'/'house/house014.jpg'/'house'
$ws.Range("A42").Value = "schultern"

# row 43: 'stürmen'/'house/house028.jpg'/'house' -> 'tauschen'/'dog/dog000.jpg'/'dog'
$ws.Range("A43").Value = "tauschen"
$ws.Range("B43").Value = "dog/dog000.jpg"
$ws.Range("C43").Value = "dog"

# row 45: 'herrschen'/'house/house007.jpg'/'house' -> 'quälen'/'house/house011.jpg'/'house'
$ws.Range("A45").Value = "quälen"
$ws.Range("B45").Value = "house/house011.jpg"

# row 46: 'quälen'/'house/house012.jpg'/'house' -> 'deuten'/'dog/dog031.jpg'/'dog'
$ws.Range("A46").Value = "deuten"
$ws.Range("B46").Value = "dog/dog031.jpg"
$ws.Range("C46").Value = "dog"

# row 48: 'packen'/'house/house016.jpg'/'house' -> 'danken'/'dog/dog002.jpg'/'dog'
$ws.Range("A48").Value = "danken"
$ws.Range("B48").Value = "dog/dog002.jpg"
$ws.Range("C48").Value = "dog"

# row 49: 'ärgern'/'house/house013.jpg'/'house' -> 'ärgern'/'house/house022.jpg'/'house'
$ws.Range("B49").Value = "house/house022.jpg"
